$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking text (e.g. "1.004", "216.69") that must
# stay plain text, matching the source inlineStr cells. Force Text number format
# before assigning so Excel does not auto-convert the string to a number, then
# drop back to the Normal style so no stray formatting is left on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.193.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.659.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5141"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.94%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2644"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.48%  "
$ws.Range("E9").Value = "  -1.91%  "
$ws.Range("E10").Value = "  -4.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07757"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.468"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.658.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.886.90"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5441"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8099"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.197.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.38%  "
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.619"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "192.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.003"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.06%  "
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1221"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.276"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.15"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.435"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05960"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.05%  "
$ws.Range("E31").Value = "  -1.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.570"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.260"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.81%  "
$ws.Range("E34").Value = "  -5.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9664"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.48%  "
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.766"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5698"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.07%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.018"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.05%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01594"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8573"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.010.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.83%  "
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.010"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.035"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05166"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.452"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.41%  "
